$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.031369823760205022
$ws.Range("B1").Value = -0.031369823786462053

$ws.Range("A2").Value = 0.0021143459219155063
$ws.Range("B2").Value = -0.0021143459740951115

$ws.Range("A3").Value = -0.045311834407592608
$ws.Range("B3").Value = 0.045311834380818101

$ws.Range("A4").Value = -0.03103266640976023
$ws.Range("B4").Value = 0.031032666375506061

$ws.Range("A5").Value = 0.031675098139486803
$ws.Range("B5").Value = -0.031675098195483538
